$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the "Folder" column (B) to "industry" for the tail of the table
# (rows 76-95), which had been left blank.
$ws.Range("B76:B95").Value = "industry"

# Make the table header row text white (header row formatting override,
# using the theme's Light 1 color rather than a hard-coded RGB so it
# matches the "Background 1" swatch used by the table style).
$ws.Range("A1:K1").Font.ThemeColor = 2  # xlThemeColorLight1

# Move the active selection back to the top of the table and clear the
# scrolled viewport that had been left at row 82.
$ws.Range("J2").Select()
